{"js": "// Replace the two generic-type label paragraphs in the bullet list:\n//   \"Resource<Statement : Resource>\"  ->  \"Resource<Statement>\"\n//   \"Resource<Kind : Statement>\"      ->  \"Resource<Kind>\"\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"Resource<Statement : Resource>\", replace: \"Resource<Statement>\" },\n  { find: \"Resource<Kind : Statement>\", replace: \"Resource<Kind>\" }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two generic-type label paragraphs in the bullet list:\n#   \"Resource<Statement : Resource>\"  ->  \"Resource<Statement>\"\n#   \"Resource<Kind : Statement>\"      ->  \"Resource<Kind>\"\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"Resource<Statement : Resource>\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"Resource<Statement>\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceOne)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Resource<Kind : Statement>\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Resource<Kind>\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceOne)\n"}
